$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff clears the numeric "Punkte"/"Note" score data for rows 11-17
# (columns C through R), leaving the Nachname/Vorname (A/B) columns intact.
$ws.Range("C11:R17").ClearContents()
